# Auto-generated script to refresh Leve profit market-price values
# across all class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3322.9333
$ws.Range("I64").Value = 3123.0908
$ws.Range("J64").Value = 3872.5
$ws.Range("K64").Value = 3123.0908
$ws.Range("L64").Value = 3872.5
$ws.Range("M64").Value = -2875.0908
$ws.Range("N64").Value = -4368.5
$ws.Range("H67").Value = 3322.9333
$ws.Range("I67").Value = 3123.0908
$ws.Range("J67").Value = 3872.5
$ws.Range("K67").Value = 3123.0908
$ws.Range("L67").Value = 3872.5
$ws.Range("M67").Value = -2265.0908
$ws.Range("N67").Value = -5588.5
$ws.Range("H98").Value = 2927194.8
$ws.Range("I98").Value = 3616.6875
$ws.Range("K98").Value = 3616.6875
$ws.Range("M98").Value = -2118.6875
$ws.Range("H122").Value = 2927194.8
$ws.Range("I122").Value = 3616.6875
$ws.Range("K122").Value = 10850.0625
$ws.Range("M122").Value = -8400.0625
$ws.Range("H125").Value = 1034.3334
$ws.Range("I125").Value = 954.75
$ws.Range("J125").Value = 1098
$ws.Range("K125").Value = 8592.75
$ws.Range("L125").Value = 9882
$ws.Range("M125").Value = -6132.75
$ws.Range("N125").Value = -14802
$ws.Range("H134").Value = 42111.11
$ws.Range("J134").Value = 42111.11
$ws.Range("L134").Value = 42111.11
$ws.Range("N134").Value = -52251.11
$ws.Range("H135").Value = 4855.8696
$ws.Range("I135").Value = 3763.8333
$ws.Range("J135").Value = 8787.200000000001
$ws.Range("K135").Value = 33874.4997
$ws.Range("L135").Value = 79084.8
$ws.Range("M135").Value = -31339.4997
$ws.Range("N135").Value = -84154.8
$ws.Range("H137").Value = 1728.6818
$ws.Range("I137").Value = 1356.625
$ws.Range("J137").Value = 1941.2858
$ws.Range("K137").Value = 4069.875
$ws.Range("L137").Value = 5823.857400000001
$ws.Range("M137").Value = -1519.875
$ws.Range("N137").Value = -10923.8574
$ws.Range("H139").Value = 78075
$ws.Range("J139").Value = 78075
$ws.Range("L139").Value = 78075
$ws.Range("N139").Value = -88355
$ws.Range("H140").Value = 97966.664
$ws.Range("J140").Value = 97966.664
$ws.Range("L140").Value = 97966.664
$ws.Range("N140").Value = -108326.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3644.5
$ws.Range("I61").Value = 3608.3845
$ws.Range("J61").Value = 4114
$ws.Range("K61").Value = 3608.3845
$ws.Range("L61").Value = 4114
$ws.Range("M61").Value = -3396.3845
$ws.Range("N61").Value = -4538
$ws.Range("H74").Value = 948.1905
$ws.Range("I74").Value = 970.6111
$ws.Range("J74").Value = 813.6667
$ws.Range("K74").Value = 970.6111
$ws.Range("L74").Value = 813.6667
$ws.Range("M74").Value = -96.61109999999996
$ws.Range("N74").Value = -2561.6667
$ws.Range("H77").Value = 948.1905
$ws.Range("I77").Value = 970.6111
$ws.Range("J77").Value = 813.6667
$ws.Range("K77").Value = 4853.055499999999
$ws.Range("L77").Value = 4068.3335
$ws.Range("M77").Value = -485.0554999999995
$ws.Range("N77").Value = -12804.3335
$ws.Range("H88").Value = 4322.636
$ws.Range("I88").Value = 1806
$ws.Range("J88").Value = 4881.8887
$ws.Range("K88").Value = 1806
$ws.Range("L88").Value = 4881.8887
$ws.Range("M88").Value = -1400
$ws.Range("N88").Value = -5693.8887
$ws.Range("H91").Value = 4322.636
$ws.Range("I91").Value = 1806
$ws.Range("J91").Value = 4881.8887
$ws.Range("K91").Value = 1806
$ws.Range("L91").Value = 4881.8887
$ws.Range("M91").Value = -402
$ws.Range("N91").Value = -7689.8887
$ws.Range("H122").Value = 3293
$ws.Range("I122").Value = 3321.2563
$ws.Range("J122").Value = 3109.3333
$ws.Range("K122").Value = 9963.768899999999
$ws.Range("L122").Value = 9327.999899999999
$ws.Range("M122").Value = -7513.768899999999
$ws.Range("N122").Value = -14227.9999
$ws.Range("H132").Value = 2902.2727
$ws.Range("I132").Value = 2634.95
$ws.Range("J132").Value = 3313.5386
$ws.Range("K132").Value = 7904.849999999999
$ws.Range("L132").Value = 9940.6158
$ws.Range("M132").Value = -5374.849999999999
$ws.Range("N132").Value = -15000.6158
$ws.Range("H136").Value = 3644.5
$ws.Range("I136").Value = 3608.3845
$ws.Range("J136").Value = 4114
$ws.Range("K136").Value = 10825.1535
$ws.Range("L136").Value = 12342
$ws.Range("M136").Value = -8275.1535
$ws.Range("N136").Value = -17442

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 46225
$ws.Range("J81").Value = 46225
$ws.Range("L81").Value = 46225
$ws.Range("N81").Value = -48347
$ws.Range("H84").Value = 46225
$ws.Range("J84").Value = 46225
$ws.Range("L84").Value = 138675
$ws.Range("N84").Value = -149283
$ws.Range("H86").Value = 5042.222
$ws.Range("I86").Value = 1900
$ws.Range("J86").Value = 11326.667
$ws.Range("K86").Value = 1900
$ws.Range("L86").Value = 11326.667
$ws.Range("M86").Value = -777
$ws.Range("N86").Value = -13572.667
$ws.Range("H89").Value = 5042.222
$ws.Range("I89").Value = 1900
$ws.Range("J89").Value = 11326.667
$ws.Range("K89").Value = 9500
$ws.Range("L89").Value = 56633.335
$ws.Range("M89").Value = -3884
$ws.Range("N89").Value = -67865.33499999999
$ws.Range("H140").Value = 86780
$ws.Range("J140").Value = 86780
$ws.Range("L140").Value = 86780
$ws.Range("N140").Value = -97140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 42997
$ws.Range("J18").Value = 42997
$ws.Range("L18").Value = 42997
$ws.Range("N18").Value = -43457
$ws.Range("H62").Value = 3969.64
$ws.Range("I62").Value = 3811.25
$ws.Range("K62").Value = 3811.25
$ws.Range("M62").Value = -3187.25
$ws.Range("H65").Value = 3969.64
$ws.Range("I65").Value = 3811.25
$ws.Range("K65").Value = 19056.25
$ws.Range("M65").Value = -15936.25
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 13159055
$ws.Range("I131").Value = 3100
$ws.Range("J131").Value = 14706815
$ws.Range("K131").Value = 9300
$ws.Range("L131").Value = 44120445
$ws.Range("M131").Value = -4260
$ws.Range("N131").Value = -44130525

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16079219
$ws.Range("I70").Value = 25575536
$ws.Range("J70").Value = 8528.385
$ws.Range("K70").Value = 25575536
$ws.Range("L70").Value = 8528.385
$ws.Range("M70").Value = -25575266
$ws.Range("N70").Value = -9068.385
$ws.Range("H73").Value = 16079219
$ws.Range("I73").Value = 25575536
$ws.Range("J73").Value = 8528.385
$ws.Range("K73").Value = 25575536
$ws.Range("L73").Value = 8528.385
$ws.Range("M73").Value = -25574600
$ws.Range("N73").Value = -10400.385
$ws.Range("H123").Value = 13396.429
$ws.Range("J123").Value = 13396.429
$ws.Range("L123").Value = 13396.429
$ws.Range("N123").Value = -18296.429
$ws.Range("H126").Value = 11112653
$ws.Range("I126").Value = 1370.3334
$ws.Range("K126").Value = 4111.0002
$ws.Range("M126").Value = -1641.0002
$ws.Range("H132").Value = 2743.64
$ws.Range("I132").Value = 2324.1177
$ws.Range("J132").Value = 3635.125
$ws.Range("K132").Value = 6972.353099999999
$ws.Range("L132").Value = 10905.375
$ws.Range("M132").Value = -4442.353099999999
$ws.Range("N132").Value = -15965.375
$ws.Range("H139").Value = 56963
$ws.Range("J139").Value = 56963
$ws.Range("L139").Value = 56963
$ws.Range("N139").Value = -67243
$ws.Range("H140").Value = 75774.5
$ws.Range("J140").Value = 75774.5
$ws.Range("L140").Value = 75774.5
$ws.Range("N140").Value = -86134.5
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4072.7273
$ws.Range("I132").Value = 4866.6665
$ws.Range("J132").Value = 3120
$ws.Range("K132").Value = 14599.9995
$ws.Range("L132").Value = 9360
$ws.Range("M132").Value = -12069.9995
$ws.Range("N132").Value = -14420
$ws.Range("H136").Value = 2901.7358
$ws.Range("I136").Value = 1895.0667
$ws.Range("J136").Value = 4214.7827
$ws.Range("K136").Value = 5685.2001
$ws.Range("L136").Value = 12644.3481
$ws.Range("M136").Value = -3135.2001
$ws.Range("N136").Value = -17744.3481
$ws.Range("H138").Value = 61350.777
$ws.Range("J138").Value = 61350.777
$ws.Range("L138").Value = 61350.777
$ws.Range("N138").Value = -71630.777
$ws.Range("H139").Value = 79600
$ws.Range("J139").Value = 79600
$ws.Range("L139").Value = 79600
$ws.Range("N139").Value = -89880

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 31964.25
$ws.Range("J123").Value = 31964.25
$ws.Range("L123").Value = 31964.25
$ws.Range("N123").Value = -41764.25
$ws.Range("H126").Value = 2001.7142
$ws.Range("I126").Value = 1093.6
$ws.Range("K126").Value = 3280.8
$ws.Range("M126").Value = -810.7999999999997
$ws.Range("H132").Value = 23440630
$ws.Range("I132").Value = 34092200
$ws.Range("J132").Value = 7180.1
$ws.Range("K132").Value = 102276600
$ws.Range("L132").Value = 21540.3
$ws.Range("M132").Value = -102274070
$ws.Range("N132").Value = -26600.3
$ws.Range("H138").Value = 104800
$ws.Range("J138").Value = 104800
$ws.Range("L138").Value = 104800
$ws.Range("N138").Value = -115080
$ws.Range("H141").Value = 79585.836
$ws.Range("J141").Value = 79585.836
$ws.Range("L141").Value = 79585.836
$ws.Range("N141").Value = -89945.836
